$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2.631648543253113, -5.290064356919663, -1.781815700103778),
    @(0.2861763049412409, -0.7290360459013701, -0.5830237417551998),
    @(-4.221912179034594, -0.03346766901844722, 9.115482657631061),
    @(-1.58611540711694, 0.8197227957620261, 4.598585906056506),
    @(2.200832122323145, -1.895753322998213, -4.244981093213744),
    @(2.756207520217482, -7.966393658191489, -2.264977694935836),
    @(2.778578069857753, -3.013215610746667, -2.033113320094297),
    @(-1.761591400025211, 3.015217549539029, 1.020017661111206),
    @(-4.130074517575332, -0.597962478681886, -1.069513860465531),
    @(-2.031331133980304, 4.770669049610229, -3.800996664631577)
)

$startRow = 22
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
